# Appends 4 new data rows (rows 9-12) to sheet1, mirroring the existing
# data layout (Date, totalScore, posWordPercentage, negWordPercentage,
# posPhrasePercentage, negPhrasePercentage, ElapsedMs, wordCount,
# sentenceCount, posWordCount, negWordCount, positivePhraseCount,
# negativePhraseCount, Method).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Excel serial date/time values (column A), matching the style (s="1")
# already used by the existing date cells in column A.
$dates = @(
    42613.758101851854,
    42613.88585648148,
    42614.884212962963,
    42615.884699074071
)

$data = @(
    @(-4, 50, 48, 50, 49, 14588, 12421, 1867, 194, 186, 13, 13),
    @(-4, 50, 48, 50, 49, 18085, 14351, 2223, 210, 203, 15, 15),
    @(-10, 52, 46, 52, 88, 12202, 9922, 1539, 155, 138, 1, 8),
    @(-4, 51, 45, 51, 77, 11763, 10966, 1470, 167, 146, 2, 7)
)

$startRow = 9

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Range("A8").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) # xlPasteFormats

    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $row[$c]
    }

    $ws.Cells.Item($r, 14).Value = "Noun"
}
